$wb = $excel.ActiveWorkbook

# Source sheet whose layout/style we replicate for the two new sheets.
$src = $wb.Worksheets.Item("syn_data2")

function Clone-DataSheet($name) {
    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $last)
    $ws.Name = $name

    # --- Row 1 labels ---
    foreach ($addr in @("A1", "L1", "W1", "AH1")) {
        $src.Range($addr).Copy()
        $ws.Range($addr).PasteSpecial(-4104)
    }

    # --- Panels (rows 2-7), each block is 10 columns wide, separated by a blank column ---
    foreach ($panel in @("A2:J7", "L2:U7", "W2:AF7", "AH2:AQ7")) {
        $topLeft = $panel.Split(":")[0]
        $src.Range($panel).Copy()
        $ws.Range($topLeft).PasteSpecial(-4122)
        $src.Range($panel).Copy()
        $ws.Range($topLeft).PasteSpecial(-4163)
    }

    return $ws
}

$ws5 = Clone-DataSheet("T10I4D100K")
$ws6 = Clone-DataSheet("kosarak")
